# Clarify the step of getting standardized aspects and units from Semantic Arts.
# Updates the two slides (4 and 13) that describe the "validate / prepare
# client reference data" step so they instead describe providing Semantic
# Arts with existing aspects/units and getting back a standardized list.

$p = $ppt.ActivePresentation

# --- Slide 13 -------------------------------------------------------------
$s13 = $p.Slides.Item(13)

# "TextBox 3" (shape id 4): reposition/resize and replace its two paragraphs.
$shp13TextBox = $s13.Shapes.Item(7)
$shp13TextBox.Left = 345.6012
$shp13TextBox.Top = 266.8739
$shp13TextBox.Width = 476.5138
$shp13TextBox.Height = 50.89221

$tr13 = $shp13TextBox.TextFrame.TextRange
$tr13.Text = "provide Semantic Arts with a list of existing aspects and units"
$tr13.InsertAfter("`rSemantic Arts will return a list of standard aspects and units")
$tr13.Font.Size = 18

# "Rectangle 10" (shape id 11): update its label text.
$shp13Rect = $s13.Shapes.Item(9)
$shp13Rect.TextFrame.TextRange.Text = "Get Semantic Arts reference data"

# --- Slide 4 ----------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# "TextBox 18" (shape id 19): widen it and replace its two paragraphs.
$shp4TextBox = $s4.Shapes.Item(5)
$shp4TextBox.Width = 364.4901

$tr4 = $shp4TextBox.TextFrame.TextRange
$tr4.Text = "provide Semantic Arts with a list of existing aspects and units"
$tr4.InsertAfter("`rSemantic Arts will return a list of standard aspects and units")
$tr4.Font.Size = 14

# "Rectangle 3" (shape id 4): update its label text.
$shp4Rect = $s4.Shapes.Item(8)
$shp4Rect.TextFrame.TextRange.Text = "Get Semantic Arts reference data"
